# "dev of account system" - add new error codes (501-510 range gets five
# new rows filled in: 505-510) to the `errors` sheet, and extend the
# trailing blank placeholder rows from 61 down to 63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: reuse of code 505, now "playerNotExistInMongo" -----------------
$ws.Range("A6").Value = "playerNotExistInMongo"
$ws.Range("C6").Value = "玩家不存在于mongo数据库"

# --- Row 7: code 506, "objectIsLocked" --------------------------------------
$ws.Range("A7").Value = "objectIsLocked"
$ws.Range("B7").Value = 506
$ws.Range("C7").Value = "对象被锁定"

# --- Row 8: code 507, "reLoginNeeded" ---------------------------------------
$ws.Range("A8").Value = "reLoginNeeded"
$ws.Range("B8").Value = 507
$ws.Range("C8").Value = "需要重新登录"

# --- Row 9: code 508, "playerAlreadyLogin" ----------------------------------
$ws.Range("A9").Value = "playerAlreadyLogin"
$ws.Range("B9").Value = 508
$ws.Range("C9").Value = "玩家已经登录"

# --- Row 10: code 509, "allianceNotExist" -----------------------------------
$ws.Range("A10").Value = "allianceNotExist"
$ws.Range("B10").Value = 509
$ws.Range("C10").Value = "联盟不存在"

# --- Row 11: code 510, "serverUnderMaintain" --------------------------------
$ws.Range("A11").Value = "serverUnderMaintain"
$ws.Range("B11").Value = 510
$ws.Range("C11").Value = "服务器维护中"

# --- Extend the trailing empty placeholder rows from 61 to 63 --------------
$newRows = $ws.Range("B62:B63")
$newRows.Borders.LineStyle = 1
$newRows.HorizontalAlignment = -4108
$newRows.VerticalAlignment = -4108
$ws.Rows.Item(62).RowHeight = 20
$ws.Rows.Item(63).RowHeight = 20

# --- Selection moves to B11, matching the author's last edit position ------
$ws.Range("B11").Select() | Out-Null
